$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date value into column A reusing the existing date style (s="1")
# by copying the format from the last pre-existing dated row (A278) before
# writing the numeric date serial, so no new style entry gets created.
function Set-DateCell($row, $serial) {
    $ws.Cells.Item(278, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 1).Value = $serial
}

# Row 279: 2018-09-16
Set-DateCell 279 43359
$ws.Cells.Item(279, 2).Value = "2018-09-16.jpg"
$ws.Cells.Item(279, 3).Value = "Gut, dass wir github benutzen. Es ist ein Unfall beim Hochladen von Dateien passiert"

# Row 280: 2018-09-21 (1)
Set-DateCell 280 43364
$ws.Cells.Item(280, 2).Value = "2018-09-21 1.JPG"
$ws.Cells.Item(280, 3).Value = "Das gehäuse für den höhensensor wurde gedruckt. Jetzt läuft er hoffentlich besser"

# Row 281: 2018-09-21 (2) -- description entered before the filename
Set-DateCell 281 43364
$ws.Cells.Item(281, 3).Value = "Nahaufnahme"
$ws.Cells.Item(281, 2).Value = "2018-09-21 2.JPG"

# Row 282: 2018-09-27 (1) -- description entered first, filename filled in later
Set-DateCell 282 43370
$ws.Cells.Item(282, 3).Value = "Das Oszilloskop ist weg. Jetzt muss ich wieder einen Arduino zum debuggen der Leitungen verwenden"

# Row 283: 2018-09-27 (2)
Set-DateCell 283 43370
$ws.Cells.Item(283, 2).Value = "2018-09-27 2.AVI"
$ws.Cells.Item(283, 3).Value = "Endlich dreht sich der Rotor wieder. Diesmal von einem anderen Anschluss"

# back to row 282 to fill in the filename
$ws.Cells.Item(282, 2).Value = "2018-09-27 1.JPG"

# Row 284: 2018-09-27 (3)
Set-DateCell 284 43370
$ws.Cells.Item(284, 2).Value = "2018-09-27 3.AVI"
$ws.Cells.Item(284, 3).Value = "Und alle 4. Das ist immer gar nicht so leicht, wie man sich es vorstellt"

# Row 285: 2018-09-29
Set-DateCell 285 43372
$ws.Cells.Item(285, 2).Value = "2018-09-29.pdf"
$ws.Cells.Item(285, 3).Value = "Die neusten (auch älteren) änderungen an der kontrollschaltung wurden upgedated. Außerdem muss der neue YMFC umgeschrieben werden, sodass ich nach 8h tippen sagen kann, dass meine finger bluten! :)"

$ws.Range("C281").Select()
